$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 21 (this pushes the old rows 21-22 down to 24-25)
$ws.Rows("21:23").Insert()

# Column A for the three new rows (same source file as the rows above)
$ws.Range("A21").Value = "tclust_wrapper1.c   "
$ws.Range("A22").Value = "tclust_wrapper1.c   "
$ws.Range("A23").Value = "tclust_wrapper1.c   "

# Fill in the new rows with data, in the order originally authored
$ws.Range("B21").Value = "double detpar"
$ws.Range("C22").Value = "double rotpar=0;"
$ws.Range("C21").Value = "double detpar=0;"
$ws.Range("B22").Value = "double rotpar"
$ws.Range("B23").Value = "double shapepar"
$ws.Range("C23").Value = "double shapepar=0;"

# Update selection to match diff (activeCell C23)
$ws.Range("C23").Select()
